$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '27.408.27'
$ws.Range('E2').Value = '  +4.50%  '
$ws.Range('D3').Value = '1.718.15'
$ws.Range('E3').Value = '  +3.77%  '
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '229.02'
$ws.Range('E5').Value = '  +4.41%  '
$ws.Range('D6').Value = '0.5377'
$ws.Range('E6').Value = '  +2.82%  '
$ws.Range('D7').Value = '1.003'
$ws.Range('E7').Value = '  -0.15%  '
$ws.Range('D8').Value = '0.2748'
$ws.Range('E8').Value = '  +3.21%  '
$ws.Range('D9').Value = '0.06709'
$ws.Range('E9').Value = '  +5.80%  '
$ws.Range('D10').Value = '21.46'
$ws.Range('E10').Value = '  +3.96%  '
$ws.Range('D11').Value = '0.07765'
$ws.Range('D12').Value = '4.711'
$ws.Range('E12').Value = '  +3.54%  '
$ws.Range('D13').Value = '1.723.54'
$ws.Range('E13').Value = '  +4.16%  '
$ws.Range('D14').Value = '1.956.21'
$ws.Range('E14').Value = '  +3.84%  '
$ws.Range('D15').Value = '0.5963'
$ws.Range('E15').Value = '  +5.07%  '
$ws.Range('D16').Value = '0.0₅8356'
$ws.Range('E16').Value = '  +3.00%  '
$ws.Range('D17').Value = '68.52'
$ws.Range('E17').Value = '  +4.66%  '
$ws.Range('D18').Value = '27.407.48'
$ws.Range('E18').Value = '  +4.58%  '
$ws.Range('B19').Value = 'Dai'
$ws.Range('C19').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D19').Value = '1.001'
$ws.Range('E19').Value = '  -0.30%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').Value = '4.792'
$ws.Range('E20').Value = '  +1.59%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').Value = '208.46'
$ws.Range('E21').Value = '  +8.34%  '
$ws.Range('D22').Value = '10.88'
$ws.Range('E22').Value = '  +5.36%  '
$ws.Range('D23').Value = '6.214'
$ws.Range('E23').Value = '  +3.03%  '
$ws.Range('D24').Value = '1.004'
$ws.Range('E24').Value = '  -0.09%  '
$ws.Range('D25').Value = '146.91'
$ws.Range('E25').Value = '  +2.31%  '
$ws.Range('D26').Value = '0.1247'
$ws.Range('E26').Value = '  +3.93%  '
$ws.Range('D27').Value = '7.394'
$ws.Range('E27').Value = '  +1.70%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').Value = '1.632'
$ws.Range('E28').Value = '  +9.15%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').Value = '16.81'
$ws.Range('E29').Value = '  +5.15%  '
$ws.Range('D30').Value = '0.05596'
$ws.Range('E30').Value = '  -0.20%  '
$ws.Range('D31').Value = '1.309'
$ws.Range('E31').Value = '  +2.46%  '
$ws.Range('D32').Value = '3.631'
$ws.Range('E32').Value = '  +3.62%  '
$ws.Range('D33').Value = '3.501'
$ws.Range('E33').Value = '  +3.60%  '
$ws.Range('D34').Value = '1.627'
$ws.Range('E34').Value = '  +2.79%  '
$ws.Range('D35').Value = '0.9693'
$ws.Range('E35').Value = '  +2.49%  '
$ws.Range('D36').Value = '2.838'
$ws.Range('E36').Value = '  +1.29%  '
$ws.Range('D37').Value = '2.434'
$ws.Range('E37').Value = '  +1.35%  '
$ws.Range('D38').Value = '0.5822'
$ws.Range('E38').Value = '  +0.83%  '
$ws.Range('D39').Value = '0.01632'
$ws.Range('E39').Value = '  +2.10%  '
$ws.Range('D40').Value = '5.869'
$ws.Range('E40').Value = '  -0.57%  '
$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').Value = '1.002'
$ws.Range('E41').Value = '  -0.16%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '1.040.99'
$ws.Range('E42').Value = '  +0.82%  '
$ws.Range('D43').Value = '0.8371'
$ws.Range('E43').Value = '  -1.03%  '
$ws.Range('D44').Value = '102.01'
$ws.Range('E44').Value = '  -0.17%  '
$ws.Range('D45').Value = '1.861.54'
$ws.Range('E45').Value = '  +3.71%  '
$ws.Range('D46').Value = '59.66'
$ws.Range('E46').Value = '  +2.09%  '
$ws.Range('D47').Value = '0.0₈107'
$ws.Range('E47').Value = '  +3.31%  '
$ws.Range('D48').Value = '8.100'
$ws.Range('E48').Value = '  +0.73%  '
$ws.Range('D49').Value = '0.9966'
$ws.Range('E49').Value = '  -0.36%  '
$ws.Range('D50').Value = '0.4367'
$ws.Range('E50').Value = '  +0.30%  '
$ws.Range('D51').Value = '0.05269'
$ws.Range('E51').Value = '  -0.87%  '
